# Fill in the teacher's info / computed "in words" total on the exam bill
# form, and tweak a couple of layout details (column A width, row 36
# height) to match the final, mostly-completed version of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header block: name / designation / department / year / term -------
$ws.Range("A3").Value = "নাম: Dr. Pintu Chandra Shill"
$ws.Range("A4").Value = "পদবী: অধ্যাপক"
$ws.Range("G4").Value = "৪র্থ"
$ws.Range("I4").Value = "১ম"
$ws.Range("B5").Value = "সিএসই"
$ws.Range("F5").Value = "বিভাগ :সিএসই"

# --- Total amount in words, under the grand-total row -------------------
$ws.Range("A32").Value = "কথায়:তের হাজার পাঁচশো সাতান্ন টাকা মাত্র।"

# --- Layout tweaks --------------------------------------------------
$ws.Columns("A").ColumnWidth = 13.43
$ws.Rows(36).RowHeight = 68.4

# --- Final selection, matching where the author left off -----------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7 | Out-Null
$ws.Range("I32").Select() | Out-Null
